$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 24,12

$data[0,0] = 1.004298570496076
$data[0,1] = 0.219705803957936
$data[0,2] = 0.07560893632798127
$data[0,3] = 0.02938583904253911
$data[0,4] = 0
$data[0,5] = 0.002673925070024397
$data[0,6] = 0
$data[0,7] = 6.044101232948748
$data[0,8] = 0
$data[0,9] = 1.112576345822077
$data[0,10] = 0.268027249794244
$data[0,11] = 0.2707570903883685

$data[1,0] = 1.007702868809133
$data[1,1] = 0.2005545605821055
$data[1,2] = 0.06903916425186196
$data[1,3] = 0.02952661281072405
$data[1,4] = 0
$data[1,5] = 0.002679789338005902
$data[1,6] = 0
$data[1,7] = 5.714031392083626
$data[1,8] = 0
$data[1,9] = 1.093446645601432
$data[1,10] = 0.2632616166814188
$data[1,11] = 0.2694741372510947

$data[2,0] = 1.01096208687477
$data[2,1] = 0.1888529921253337
$data[2,2] = 0.0650484898235959
$data[2,3] = 0.02961878800619672
$data[2,4] = 0
$data[2,5] = 0.002683571347395843
$data[2,6] = 0
$data[2,7] = 5.51068234846889
$data[2,8] = 0
$data[2,9] = 1.082774347797141
$data[2,10] = 0.2604413552815217
$data[2,11] = 0.2689148003331674

$data[3,0] = 1.012583182297988
$data[3,1] = 0.184098545961632
$data[3,2] = 0.06343294824789325
$data[3,3] = 0.02965779748381725
$data[3,4] = 0
$data[3,5] = 0.002685158317402206
$data[3,6] = 0
$data[3,7] = 5.427634953416145
$data[3,8] = 0
$data[3,9] = 1.078694008875061
$data[3,10] = 0.259318629671327
$data[3,11] = 0.2687441646320288

$data[4,0] = 1.012870027614667
$data[4,1] = 0.1833099104393341
$data[4,2] = 0.06316532956887499
$data[4,3] = 0.02966436251408072
$data[4,4] = 0
$data[4,5] = 0.002685424601800401
$data[4,6] = 0
$data[4,7] = 5.413833809164572
$data[4,8] = 0
$data[4,9] = 1.078032666230996
$data[4,10] = 0.2591338038970505
$data[4,11] = 0.2687192875469222

$data[5,0] = 1.010982764783165
$data[5,1] = 0.1887888156649069
$data[5,2] = 0.06502665899917304
$data[5,3] = 0.0296193082357416
$data[5,4] = 0
$data[5,5] = 0.002683592564151541
$data[5,6] = 0
$data[5,7] = 5.509563086721386
$data[5,8] = 0
$data[5,9] = 1.082718232521842
$data[5,10] = 0.2604261063768618
$data[5,11] = 0.2689122672589015

$data[6,0] = 1.005229136023132
$data[6,1] = 0.2130903060053981
$data[6,2] = 0.07333460899845079
$data[6,3] = 0.0294331890493173
$data[6,4] = 0
$data[6,5] = 0.002675909540963902
$data[6,6] = 0
$data[6,7] = 5.930429869587641
$data[6,8] = 0
$data[6,9] = 1.105756891236382
$data[6,10] = 0.2663620472332013
$data[6,11] = 0.2702672395142685

$data[7,0] = 1.003270041043294
$data[7,1] = 0.2612208456589542
$data[7,2] = 0.08997783168203455
$data[7,3] = 0.02911356822302785
$data[7,4] = 0
$data[7,5] = 0.002662273881923447
$data[7,6] = 0
$data[7,7] = 6.750751370329482
$data[7,8] = 0
$data[7,9] = 1.159513009831841
$data[7,10] = 0.278846222161377
$data[7,11] = 0.2747435567707797

$data[8,0] = 1.007586410388342
$data[8,1] = 0.2969038916173758
$data[8,2] = 0.102433171083419
$data[8,3] = 0.02890614522330259
$data[8,4] = 0
$data[8,5] = 0.002653116707248984
$data[8,6] = 0
$data[8,7] = 7.351057501380865
$data[8,8] = 0
$data[8,9] = 1.204328479752121
$data[8,10] = 0.2885395468487673
$data[8,11] = 0.2791522669366202

$data[9,0] = 1.01081605970839
$data[9,1] = 0.3132140371980086
$data[9,2] = 0.1081518564441382
$data[9,3] = 0.02881768091288617
$data[9,4] = 0
$data[9,5] = 0.002649135399668317
$data[9,6] = 0
$data[9,7] = 7.623782787143057
$data[9,8] = 0
$data[9,9] = 1.225892425857438
$data[9,10] = 0.2930640683485279
$data[9,11] = 0.2814036139702836

$data[10,0] = 1.012222427600221
$data[10,1] = 0.31940196416096
$data[10,2] = 0.1103251733168946
$data[10,3] = 0.02878502521058401
$data[10,4] = 0
$data[10,5] = 0.002647654103313585
$data[10,6] = 0
$data[10,7] = 7.727016495331668
$data[10,8] = 0
$data[10,9] = 1.234229023578933
$data[10,10] = 0.2947940304423895
$data[10,11] = 0.2822916818910457

$data[11,0] = 1.011911364546677
$data[11,1] = 0.3180687589358797
$data[11,2] = 0.1098567616565873
$data[11,3] = 0.02879202072951914
$data[11,4] = 0
$data[11,5] = 0.002647971958321712
$data[11,6] = 0
$data[11,7] = 7.704784952756427
$data[11,8] = 0
$data[11,9] = 1.232425967048528
$data[11,10] = 0.2944207113219335
$data[11,11] = 0.2820988372734732

$data[12,0] = 1.010928079844746
$data[12,1] = 0.3137228860098844
$data[12,2] = 0.1083304995614327
$data[12,3] = 0.02881497741823524
$data[12,4] = 0
$data[12,5] = 0.002649013005413355
$data[12,6] = 0
$data[12,7] = 7.632276659636887
$data[12,8] = 0
$data[12,9] = 1.226574850822544
$data[12,10] = 0.2932060596716042
$data[12,11] = 0.281475962660636

$data[13,0] = 1.010349709029128
$data[13,1] = 0.3110624412537675
$data[13,2] = 0.1073966384260672
$data[13,3] = 0.02882914883435528
$data[13,4] = 0
$data[13,5] = 0.002649654102137734
$data[13,6] = 0
$data[13,7] = 7.587858136587158
$data[13,8] = 0
$data[13,9] = 1.223013162225783
$data[13,10] = 0.2924642181374537
$data[13,11] = 0.2810990668578341

$data[14,0] = 1.007400933343092
$data[14,1] = 0.2958395878848705
$data[14,2] = 0.1020605186368471
$data[14,3] = 0.02891204485156751
$data[14,4] = 0
$data[14,5] = 0.002653380590177565
$data[14,6] = 0
$data[14,7] = 7.333227921072364
$data[14,8] = 0
$data[14,9] = 1.202943042132915
$data[14,10] = 0.2882461800250837
$data[14,11] = 0.2790100980108576

$data[15,0] = 1.005917153876055
$data[15,1] = 0.2865210899685167
$data[15,2] = 0.09880061665732853
$data[15,3] = 0.02896440571017844
$data[15,4] = 0
$data[15,5] = 0.002655713762590905
$data[15,6] = 0
$data[15,7] = 7.17693608544775
$data[15,8] = 0
$data[15,9] = 1.190933245754735
$data[15,10] = 0.2856880621841924
$data[15,11] = 0.2777916693089324

$data[16,0] = 1.005182815510864
$data[16,1] = 0.2811686271248277
$data[16,2] = 0.0969305544965664
$data[16,3] = 0.02899507722708128
$data[16,4] = 0
$data[16,5] = 0.002657073102254915
$data[16,6] = 0
$data[16,7] = 7.087006952091656
$data[16,8] = 0
$data[16,9] = 1.184136235228806
$data[16,10] = 0.28422751725914
$data[16,11] = 0.277113985991182

$data[17,0] = 1.004954596168858
$data[17,1] = 0.2793576138498111
$data[17,2] = 0.09629822730180138
$data[17,3] = 0.02900555749856049
$data[17,4] = 0
$data[17,5] = 0.002657536338334018
$data[17,6] = 0
$data[17,7] = 7.056552319620408
$data[17,8] = 0
$data[17,9] = 1.181853852415884
$data[17,10] = 0.2837348569550642
$data[17,11] = 0.2768885001883703

$data[18,0] = 1.006062769790702
$data[18,1] = 0.2875123022006676
$data[18,2] = 0.09914712534481396
$data[18,3] = 0.02895877440030858
$data[18,4] = 0
$data[18,5] = 0.002655463596763371
$data[18,6] = 0
$data[18,7] = 7.1935770869338
$data[18,8] = 0
$data[18,9] = 1.192200240443213
$data[18,10] = 0.2859592580572468
$data[18,11] = 0.2779189784876124

$data[19,0] = 1.011211907225771
$data[19,1] = 0.3149990550083999
$data[19,2] = 0.1087785871294784
$data[19,3] = 0.02880821161040803
$data[19,4] = 0
$data[19,5] = 0.002648706511122327
$data[19,6] = 0
$data[19,7] = 7.653575147861488
$data[19,8] = 0
$data[19,9] = 1.228288817698086
$data[19,10] = 0.2935623803977592
$data[19,11] = 0.2816579502045542

$data[20,0] = 1.01564661558416
$data[20,1] = 0.3330312336319707
$data[20,2] = 0.1151187168317875
$data[20,3] = 0.02871472683937182
$data[20,4] = 0
$data[20,5] = 0.002644443808189471
$data[20,6] = 0
$data[20,7] = 7.95397365325033
$data[20,8] = 0
$data[20,9] = 1.252871074402719
$data[20,10] = 0.2986284022595811
$data[20,11] = 0.2843087806097842

$data[21,0] = 1.0131814335858
$data[21,1] = 0.3234007480951391
$data[21,2] = 0.111730648461787
$data[21,3] = 0.02876417273881793
$data[21,4] = 0
$data[21,5] = 0.002646704908679577
$data[21,6] = 0
$data[21,7] = 7.79366367150601
$data[21,8] = 0
$data[21,9] = 1.239659387048647
$data[21,10] = 0.2959156701090961
$data[21,11] = 0.2828749611392354

$data[22,0] = 1.005996567190891
$data[22,1] = 0.2870641600279953
$data[22,2] = 0.09899045601473233
$data[22,3] = 0.02896131854435469
$data[22,4] = 0
$data[22,5] = 0.002655576640719914
$data[22,6] = 0
$data[22,7] = 7.186053924126782
$data[22,8] = 0
$data[22,9] = 1.1916270971671
$data[22,10] = 0.2858366187916346
$data[22,11] = 0.2778613509814321

$data[23,0] = 1.002794558178664
$data[23,1] = 0.2481461182508724
$data[23,2] = 0.08543641570115312
$data[23,3] = 0.02919520412556487
$data[23,4] = 0
$data[23,5] = 0.002665810685307974
$data[23,6] = 0
$data[23,7] = 6.529299756977366
$data[23,8] = 0
$data[23,9] = 1.144042926915091
$data[23,10] = 0.2753778970628389
$data[23,11] = 0.2733367269515625

$ws.Range("B2:M25").Value = $data
